$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values for rows 3-6, and column B for row 3
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2

# Row 7 no longer exists in the new data range - delete it (shifts nothing below, just clears it)
$ws.Rows("7").Delete()
